$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old third row entirely (old layout used rows 1-3, new layout only needs 1-2)
$ws.Rows(3).Delete()

# ---------- Row 1 : headers ----------
# A1 already carries the bold/filled/bordered header style - just change its text
$ws.Range("A1").Value = "testcase"
# B1 already carries the same header style - just change its text
$ws.Range("B1").Value = "RGID"

# C1:F1 are brand-new header cells - clone A1's header formatting onto them
$ws.Range("A1").Copy()
$ws.Range("C1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C1").Value = "ProviderID"
$ws.Range("D1").Value = "Count"
$ws.Range("E1").Value = "Run "
$ws.Range("F1").Value = "Environment"
# B1:F1 are centered headers
$ws.Range("B1:F1").HorizontalAlignment = -4108

# G1 - same (non-centered) header look as A1, same text as A1
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G1").Value = "testcase"

# ---------- Row 2 : data ----------
# B2:F2 clone the old bordered data-cell look (A2 already carries it), then get centered
$ws.Range("A2").Copy()
$ws.Range("B2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A2").ClearContents()
$ws.Range("B2").Value = "RG-35022"
$ws.Range("C2").Value = "P-157855884725"
$ws.Range("D2").Value = """3"""
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "QA"
$ws.Range("B2:F2").HorizontalAlignment = -4108

$ws.Range("D3").Select()
